# Scen_PV_CAP_BND.xlsx - "Add files via upload" edit
#
# 1) Rename the process NEW_BC_PP -> NEW_WINDOFF_PP (cell B10 on Sheet1,
#    which is the sole use of that shared string).
# 2) Update the saved view state: zoom 265% -> 220%, and move the
#    active selection from C12 to D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- content change -------------------------------------------------
$ws.Range("B10").Value = "NEW_WINDOFF_PP"

# --- view/selection change ------------------------------------------
$null = $ws.Activate()
$excel.ActiveWindow.Zoom = 220
$null = $ws.Range("D12").Select()
